# Apply the changes described by the diff:
# 1) Sheet "Resumen": C2 -> new Maximo value
# 2) Sheet "Solucion": reorder/update the "Salida" (column B) values for rows 2-33
# 3) Sheet "Metricas": B2 -> new Z1 time, B3 -> new Z2 time

$wb = $excel.ActiveWorkbook

# --- Sheet: Resumen ---
$wsResumen = $wb.Worksheets.Item("Resumen")
$wsResumen.Range("C2").Value = 562.2322953832843

# --- Sheet: Solucion ---
$wsSolucion = $wb.Worksheets.Item("Solucion")

$salidaValues = @{
    2  = "S031"
    3  = "S036"
    5  = "S039"
    6  = "S034"
    7  = "S026"
    8  = "S035"
    9  = "S040"
    11 = "S006"
    12 = "S030"
    14 = "S028"
    15 = "S003"
    17 = "S002"
    18 = "S025"
    19 = "S004"
    20 = "S033"
    22 = "S038"
    23 = "S001"
    24 = "S008"
    25 = "S007"
    26 = "S009"
    27 = "S013"
    28 = "S014"
    29 = "S011"
    30 = "S016"
    32 = "S010"
    33 = "S015"
}

foreach ($row in $salidaValues.Keys) {
    $wsSolucion.Cells.Item($row, 2).Value = $salidaValues[$row]
}

# --- Sheet: Metricas ---
$wsMetricas = $wb.Worksheets.Item("Metricas")
$wsMetricas.Range("B2").Value = 562.2322953832843
$wsMetricas.Range("B3").Value = 542.9935939020434
